$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 8197.77
$ws.Range("I41").Value = 1004.2
$ws.Range("J41").Value = 12693.75
$ws.Range("K41").Value = 1004.2
$ws.Range("L41").Value = 12693.75
$ws.Range("M41").Value = -564.2
$ws.Range("N41").Value = -13573.75
$ws.Range("H62").Value = 2795.75
$ws.Range("I62").Value = 2273.2
$ws.Range("J62").Value = 3666.6667
$ws.Range("K62").Value = 2273.2
$ws.Range("L62").Value = 3666.6667
$ws.Range("M62").Value = -1649.2
$ws.Range("N62").Value = -4914.6667
$ws.Range("H65").Value = 2795.75
$ws.Range("I65").Value = 2273.2
$ws.Range("J65").Value = 3666.6667
$ws.Range("K65").Value = 11366
$ws.Range("L65").Value = 18333.3335
$ws.Range("M65").Value = -8246
$ws.Range("N65").Value = -24573.3335
$ws.Range("H103").Value = 942.9091
$ws.Range("J103").Value = 1046
$ws.Range("L103").Value = 3138
$ws.Range("N103").Value = -4310
$ws.Range("H116").Value = 11667
$ws.Range("I116").Value = 11001.5
$ws.Range("K116").Value = 11001.5
$ws.Range("M116").Value = -7559.5
$ws.Range("H132").Value = 11405.509
$ws.Range("I132").Value = 1903.75
$ws.Range("K132").Value = 5711.25
$ws.Range("M132").Value = -3181.25
$ws.Range("H137").Value = 2584.4517
$ws.Range("I137").Value = 2312.2693
$ws.Range("J137").Value = 3999.8
$ws.Range("K137").Value = 6936.8079
$ws.Range("L137").Value = 11999.4
$ws.Range("M137").Value = -4386.8079
$ws.Range("N137").Value = -17099.4
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4848.8887
$ws.Range("I45").Value = 1828.2
$ws.Range("K45").Value = 1828.2
$ws.Range("M45").Value = -1451.2
$ws.Range("H61").Value = 1934.8823
$ws.Range("I61").Value = 1934.8823
$ws.Range("K61").Value = 1934.8823
$ws.Range("M61").Value = -1722.8823
$ws.Range("H110").Value = 4083.6365
$ws.Range("I110").Value = 3966.7222
$ws.Range("K110").Value = 3966.7222
$ws.Range("M110").Value = -1921.7222
$ws.Range("H118").Value = 49800.31
$ws.Range("J118").Value = 49800.31
$ws.Range("L118").Value = 49800.31
$ws.Range("N118").Value = -53114.31
$ws.Range("H122").Value = 5205.136
$ws.Range("I122").Value = 5608.8887
$ws.Range("J122").Value = 4925.615
$ws.Range("K122").Value = 16826.6661
$ws.Range("L122").Value = 14776.845
$ws.Range("M122").Value = -14376.6661
$ws.Range("N122").Value = -19676.845
$ws.Range("H132").Value = 1613.9286
$ws.Range("I132").Value = 1667.8948
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 5003.6844
$ws.Range("L132").Value = 4500
$ws.Range("M132").Value = -2473.6844
$ws.Range("N132").Value = -9560
$ws.Range("H136").Value = 1934.8823
$ws.Range("I136").Value = 1934.8823
$ws.Range("K136").Value = 5804.6469
$ws.Range("M136").Value = -3254.6469
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2367.9412
$ws.Range("I20").Value = 1534.6364
$ws.Range("K20").Value = 1534.6364
$ws.Range("M20").Value = -1287.6364
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H94").Value = 3796.1538
$ws.Range("I94").Value = 3796.1538
$ws.Range("K94").Value = 3796.1538
$ws.Range("M94").Value = -3345.1538
$ws.Range("H107").Value = 6151.55
$ws.Range("I107").Value = 5538.75
$ws.Range("K107").Value = 5538.75
$ws.Range("M107").Value = -3618.75
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2593.7827
$ws.Range("I16").Value = 2376.842
$ws.Range("K16").Value = 2376.842
$ws.Range("M16").Value = -2089.842
$ws.Range("H31").Value = 1253.0238
$ws.Range("J31").Value = 1600.4445
$ws.Range("L31").Value = 1600.4445
$ws.Range("N31").Value = -2190.4445
$ws.Range("H34").Value = 1253.0238
$ws.Range("J34").Value = 1600.4445
$ws.Range("L34").Value = 1600.4445
$ws.Range("N34").Value = -2004.4445
$ws.Range("H105").Value = 3000
$ws.Range("I105").Value = 3000
$ws.Range("K105").Value = 3000
$ws.Range("M105").Value = -1253
$ws.Range("H107").Value = 5423
$ws.Range("I107").Value = 1041.6
$ws.Range("J107").Value = 9074.166999999999
$ws.Range("K107").Value = 1041.6
$ws.Range("L107").Value = 9074.166999999999
$ws.Range("M107").Value = 878.4000000000001
$ws.Range("N107").Value = -12914.167
$ws.Range("H113").Value = 2593.7827
$ws.Range("I113").Value = 2376.842
$ws.Range("K113").Value = 2376.842
$ws.Range("M113").Value = -206.8420000000001
$ws.Range("H132").Value = 2983.2964
$ws.Range("I132").Value = 2269.4119
$ws.Range("K132").Value = 6808.2357
$ws.Range("M132").Value = -4278.2357
$ws.Range("H134").Value = 3907.6316
$ws.Range("I134").Value = 3822.4412
$ws.Range("J134").Value = 4631.75
$ws.Range("K134").Value = 11467.3236
$ws.Range("L134").Value = 13895.25
$ws.Range("M134").Value = -8932.3236
$ws.Range("N134").Value = -18965.25
$ws.Range("H135").Value = 70777.7
$ws.Range("J135").Value = 70777.7
$ws.Range("L135").Value = 70777.7
$ws.Range("N135").Value = -80917.7
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H95").Value = 10000
$ws.Range("J95").Value = 10000
$ws.Range("L95").Value = 30000
$ws.Range("N95").Value = -34118
$ws.Range("H139").Value = 500004500
$ws.Range("I139").Value = 1000000000
$ws.Range("J139").Value = 9000
$ws.Range("K139").Value = 3000000000
$ws.Range("L139").Value = 27000
$ws.Range("M139").Value = -2999994860
$ws.Range("N139").Value = -37280
$ws.Range("H141").Value = 111116630
$ws.Range("J141").Value = 10000
$ws.Range("L141").Value = 30000
$ws.Range("N141").Value = -40360
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 144213.75
$ws.Range("I70").Value = 190109.17
$ws.Range("J70").Value = 6527.5
$ws.Range("K70").Value = 190109.17
$ws.Range("L70").Value = 6527.5
$ws.Range("M70").Value = -189839.17
$ws.Range("N70").Value = -7067.5
$ws.Range("H73").Value = 144213.75
$ws.Range("I73").Value = 190109.17
$ws.Range("J73").Value = 6527.5
$ws.Range("K73").Value = 190109.17
$ws.Range("L73").Value = 6527.5
$ws.Range("M73").Value = -189173.17
$ws.Range("N73").Value = -8399.5
$ws.Range("H113").Value = 10155.556
$ws.Range("I113").Value = 4280
$ws.Range("J113").Value = 17500
$ws.Range("K113").Value = 4280
$ws.Range("L113").Value = 17500
$ws.Range("M113").Value = -2110
$ws.Range("N113").Value = -21840
$ws.Range("H122").Value = 4283.5386
$ws.Range("I122").Value = 3768.7
$ws.Range("J122").Value = 5999.6665
$ws.Range("K122").Value = 11306.1
$ws.Range("L122").Value = 17998.9995
$ws.Range("M122").Value = -8856.099999999999
$ws.Range("N122").Value = -22898.9995
$ws.Range("H123").Value = 39496.25
$ws.Range("J123").Value = 39496.25
$ws.Range("L123").Value = 39496.25
$ws.Range("N123").Value = -44396.25
$ws.Range("H126").Value = 4653.143
$ws.Range("I126").Value = 4653.143
$ws.Range("K126").Value = 13959.429
$ws.Range("M126").Value = -11489.429
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7041.6665
$ws.Range("I122").Value = 5125
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 15375
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -12925
$ws.Range("N122").Value = -28900
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 42659.69
$ws.Range("I96").Value = 128489.75
$ws.Range("J96").Value = 4513
$ws.Range("K96").Value = 128489.75
$ws.Range("L96").Value = 4513
$ws.Range("M96").Value = -127116.75
$ws.Range("N96").Value = -7259
$ws.Range("H100").Value = 1189.9
$ws.Range("I100").Value = 1630.2
$ws.Range("J100").Value = 749.6
$ws.Range("K100").Value = 3260.4
$ws.Range("L100").Value = 1499.2
$ws.Range("M100").Value = -2719.4
$ws.Range("N100").Value = -2581.2
$ws.Range("H136").Value = 805.3333
$ws.Range("I136").Value = 649.04346
$ws.Range("J136").Value = 4400
$ws.Range("K136").Value = 1947.13038
$ws.Range("L136").Value = 13200
$ws.Range("M136").Value = 602.8696199999999
$ws.Range("N136").Value = -18300
